# Generate Report for Handback
#
# The localization-status workbook tracks, per language sheet, the
# handoff of source files and (new in this edit) the handback of the
# translated files. This script:
#   1. Updates the "Status" text (shared by every row on every sheet)
#      from "Ready for handoff" to "Handed back: in sync with en-US".
#   2. Populates the new "Latest Target File" (F) and "Latest Handback
#      File" (G) columns for every data row, each as a hyperlinked file
#      name (mirroring the existing A/B/D hyperlink columns).
#   3. Fills in the "Latest Handback DateTime" column (H) now that the
#      files have actually been handed back.

$wb = $excel.ActiveWorkbook

# --- 1. Status text, shared by every row on every sheet (Overview included) ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

function Set-HandbackRow($SheetName, $Row, $TargetFile, $TargetUrl, $HandbackFile, $HandbackUrl, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # F: Latest Target File (the .md source, now also the handback target)
    $fCell = $ws.Range("F$Row")
    $fCell.Value = $TargetFile
    $ws.Hyperlinks.Add($fCell, $TargetUrl, "", "", $TargetFile) | Out-Null

    # G: Latest Handback File (the localized .xlf that was handed back)
    $gCell = $ws.Range("G$Row")
    $gCell.Value = $HandbackFile
    $ws.Hyperlinks.Add($gCell, $HandbackUrl, "", "", $HandbackFile) | Out-Null

    # H: Latest Handback DateTime
    $ws.Range("H$Row").Value = $HandbackDateTime
}

# --- zh-cn sheet ---
Set-HandbackRow "zh-cn" 2 `
    "15d11c1e-0815-4114-bde9-2a6d240f3db1.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8c864bce8b1f69082422c87895c1e139b25b9567/e2e/15d11c1e-0815-4114-bde9-2a6d240f3db1.md" `
    "15d11c1e-0815-4114-bde9-2a6d240f3db1.271158a45f2da6e9e17724a7ccc876ddeb02ef33.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6f1cee89f1f8e81ca14b1b8fb3de9a7e344877fb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/15d11c1e-0815-4114-bde9-2a6d240f3db1.271158a45f2da6e9e17724a7ccc876ddeb02ef33.zh-cn.xlf" `
    "2016-03-19 16:14:43"

Set-HandbackRow "zh-cn" 3 `
    "60820c3e-091d-474e-924c-ec14ccb9167a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8c864bce8b1f69082422c87895c1e139b25b9567/e2e/60820c3e-091d-474e-924c-ec14ccb9167a.md" `
    "60820c3e-091d-474e-924c-ec14ccb9167a.eb26a0fa9fc2147fa9f73846944daca0d905b35c.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6f1cee89f1f8e81ca14b1b8fb3de9a7e344877fb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/60820c3e-091d-474e-924c-ec14ccb9167a.eb26a0fa9fc2147fa9f73846944daca0d905b35c.zh-cn.xlf" `
    "2016-03-19 16:14:43"

# --- de-de sheet ---
Set-HandbackRow "de-de" 2 `
    "15d11c1e-0815-4114-bde9-2a6d240f3db1.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8c864bce8b1f69082422c87895c1e139b25b9567/e2e/15d11c1e-0815-4114-bde9-2a6d240f3db1.md" `
    "15d11c1e-0815-4114-bde9-2a6d240f3db1.271158a45f2da6e9e17724a7ccc876ddeb02ef33.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/21b15c998a8ce9ab9a017ebab84fbf1dbd629ae8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/15d11c1e-0815-4114-bde9-2a6d240f3db1.271158a45f2da6e9e17724a7ccc876ddeb02ef33.de-de.xlf" `
    "2016-03-19 16:14:48"

Set-HandbackRow "de-de" 3 `
    "60820c3e-091d-474e-924c-ec14ccb9167a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8c864bce8b1f69082422c87895c1e139b25b9567/e2e/60820c3e-091d-474e-924c-ec14ccb9167a.md" `
    "60820c3e-091d-474e-924c-ec14ccb9167a.eb26a0fa9fc2147fa9f73846944daca0d905b35c.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/21b15c998a8ce9ab9a017ebab84fbf1dbd629ae8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/60820c3e-091d-474e-924c-ec14ccb9167a.eb26a0fa9fc2147fa9f73846944daca0d905b35c.de-de.xlf" `
    "2016-03-19 16:14:48"

Write-Host "Handback report generated."
